$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5006976723670959
$ws.Range("B1").Value = 0.5752987265586853
$ws.Range("C1").Value = 4.676632881164551
$ws.Range("D1").Value = 1.882308483123779
$ws.Range("E1").Value = 0.896948516368866
